$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$sub3 = [char]0x2083

$ws.Cells.Item(2, 4).Value = "34.336.11"
$ws.Cells.Item(2, 5).Value = "  +11.85%  "
$ws.Cells.Item(3, 4).Value = "1.819.26"
$ws.Cells.Item(3, 5).Value = "  +8.47%  "
$ws.Cells.Item(4, 4).Value = "'0.998"
$ws.Cells.Item(4, 5).Value = "  -0.09%  "
$ws.Cells.Item(5, 4).Value = "'229.42"
$ws.Cells.Item(5, 5).Value = "  +4.78%  "
$ws.Cells.Item(6, 4).Value = "'0.573"
$ws.Cells.Item(6, 5).Value = "  +8.58%  "
$ws.Cells.Item(7, 4).Value = "'0.998"
$ws.Cells.Item(7, 5).Value = "  +0.02%  "
$ws.Cells.Item(8, 4).Value = "'31.59"
$ws.Cells.Item(8, 5).Value = "  +8.65%  "
$ws.Cells.Item(9, 4).Value = "'46.67"
$ws.Cells.Item(9, 5).Value = "  +3.93%  "
$ws.Cells.Item(10, 4).Value = "'0.287"
$ws.Cells.Item(10, 5).Value = "  +8.97%  "
$ws.Cells.Item(11, 4).Value = "'0.0676"
$ws.Cells.Item(11, 5).Value = "  +5.18%  "
$ws.Cells.Item(12, 4).Value = "'0.0932"
$ws.Cells.Item(12, 5).Value = "  +3.12%  "
$ws.Cells.Item(13, 4).Value = "2.077.73"
$ws.Cells.Item(13, 5).Value = "  +8.28%  "
$ws.Cells.Item(14, 4).Value = "1.821.20"
$ws.Cells.Item(14, 5).Value = "  +8.86%  "
$ws.Cells.Item(15, 4).Value = "'0.647"
$ws.Cells.Item(15, 5).Value = "  +6.88%  "
$ws.Cells.Item(16, 4).Value = "34.253.09"
$ws.Cells.Item(17, 4).Value = "'10.25"
$ws.Cells.Item(17, 5).Value = "  +1.36%  "
$ws.Cells.Item(18, 5).Value = "  +6.39%  "
$ws.Cells.Item(19, 4).Value = "'70.43"
$ws.Cells.Item(19, 5).Value = "  +6.81%  "
$ws.Cells.Item(20, 4).Value = "'258.94"
$ws.Cells.Item(20, 5).Value = "  +6.57%  "
$ws.Cells.Item(21, 4).Value = "0.0{0}0755" -f $sub3
$ws.Cells.Item(21, 5).Value = "  +5.12%  "
$ws.Cells.Item(22, 4).Value = "'0.999"
$ws.Cells.Item(22, 5).Value = "  +0.00%  "
$ws.Cells.Item(23, 4).Value = "'10.64"
$ws.Cells.Item(23, 5).Value = "  +6.97%  "
$ws.Cells.Item(24, 4).Value = "'4.35"
$ws.Cells.Item(24, 5).Value = "  +2.85%  "
$ws.Cells.Item(25, 5).Value = "  +2.41%  "
$ws.Cells.Item(26, 4).Value = "'159.35"
$ws.Cells.Item(26, 5).Value = "  +0.20%  "
$ws.Cells.Item(27, 4).Value = "'16.77"
$ws.Cells.Item(27, 5).Value = "  +6.03%  "
$ws.Cells.Item(28, 5).Value = "  +5.37%  "
$ws.Cells.Item(29, 4).Value = "'7.07"
$ws.Cells.Item(29, 5).Value = "  +5.86%  "
$ws.Cells.Item(30, 4).Value = "'0.998"
$ws.Cells.Item(30, 5).Value = "  -0.16%  "
$ws.Cells.Item(31, 4).Value = "'3.90"
$ws.Cells.Item(31, 5).Value = "  +12.82%  "
$ws.Cells.Item(32, 5).Value = "  +6.68%  "
$ws.Cells.Item(33, 5).Value = "  +6.11%  "
$ws.Cells.Item(34, 4).Value = "'3.59"
$ws.Cells.Item(34, 5).Value = "  +8.68%  "
$ws.Cells.Item(35, 4).Value = "1.530.15"
$ws.Cells.Item(35, 5).Value = "  +0.84%  "
$ws.Cells.Item(36, 5).Value = "  +1.80%  "
$ws.Cells.Item(37, 5).Value = "  +5.53%  "
$ws.Cells.Item(38, 4).Value = "'0.636"
$ws.Cells.Item(38, 5).Value = "  +5.80%  "
$ws.Cells.Item(39, 4).Value = "'0.0190"
$ws.Cells.Item(39, 5).Value = "  +6.62%  "
$ws.Cells.Item(40, 4).Value = "'83.79"
$ws.Cells.Item(40, 5).Value = "  +0.96%  "
$ws.Cells.Item(41, 4).Value = "'2.80"
$ws.Cells.Item(41, 5).Value = "  +4.83%  "
$ws.Cells.Item(42, 5).Value = "  +3.00%  "
$ws.Cells.Item(43, 4).Value = "'0.912"
$ws.Cells.Item(43, 5).Value = "  +9.09%  "
$ws.Cells.Item(44, 5).Value = "  +5.33%  "
$ws.Cells.Item(45, 4).Value = "'0.0521"
$ws.Cells.Item(45, 5).Value = "  +4.44%  "
$ws.Cells.Item(46, 5).Value = "  +5.85%  "
$ws.Cells.Item(47, 4).Value = "1.969.09"
$ws.Cells.Item(48, 4).Value = "'5.83"
$ws.Cells.Item(48, 5).Value = "  +5.30%  "
$ws.Cells.Item(49, 4).Value = "'12.16"
$ws.Cells.Item(49, 5).Value = "  +18.12%  "
$ws.Cells.Item(50, 5).Value = "  +0.00%  "
$ws.Cells.Item(51, 4).Value = "'52.05"
$ws.Cells.Item(51, 5).Value = "  +4.33%  "

Write-Host "done"
